$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reference cell that already carries the "command name" style (italic,
# 8pt, blue Consolas) used throughout column G/K - copy it onto the new
# cells so the added rows look the same as the existing ones.
# ---------------------------------------------------------------------------
$styleSource = $ws.Range("G3")

function Set-CmdCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $styleSource.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

# New block below the existing table (rows 22-39), per the revised
# "Realizace Zakladni Deska" sequence.
$ws.Range("F22").Value = "THL Katherine"

Set-CmdCell "G23" "CMD_SET_ALL_PIXEL_CONFIG"
Set-CmdCell "G24" "CMD_HW_COMMAND_START"
Set-CmdCell "G25" "CMD_ACQ_TIME_SETTING_LSB"
Set-CmdCell "G26" "CMD_NUMBER_OF_FRAMES_SETTING"
Set-CmdCell "G27" "CMD_BIAS_SETTING"
Set-CmdCell "G28" "CMD_GET_BIAS_VOLTAGE"

$ws.Range("F29").Value = "5x"
Set-CmdCell "G29" "CMD_ACQUISITION_SETUP"

Set-CmdCell "G30" "CMD_INTERNAL_TDC_SETTINGS"
Set-CmdCell "G31" "CMD_GET_ACQUISITION_SETUP"
Set-CmdCell "G32" "CMD_HW_COMMAND_START"
Set-CmdCell "G33" "CMD_INTERNAL_DAC_SETTINGS"
Set-CmdCell "G34" "CMD_HW_COMMAND_START"
Set-CmdCell "G35" "CMD_TPX2_SET_OMR"
Set-CmdCell "G36" "CMD_ACQ_MODE_SETTING"
Set-CmdCell "G37" "CMD_TPX2_SET_FRQ"
Set-CmdCell "G38" "CMD_TPX2_SET_COL_TRIGGER"
Set-CmdCell "G39" "CMD_ACQ_START"

# Selection state, matching the saved view after the edit.
$ws.Range("L32").Select()
